$d = $word.ActiveDocument

$replacements = @(
    @("31÷9=3, 4", "61÷9=6, 7"),
    @("44÷8=5, 4", "89÷2=44, 1"),
    @("99÷6=16, 3", "80÷6=13, 2"),
    @("72÷3=24, 0", "95÷4=23, 3"),
    @("77÷8=9, 5", "90÷7=12, 6"),
    @("18÷6=3, 0", "80÷6=13, 2"),
    @("78÷8=9, 6", "98÷6=16, 2"),
    @("34÷2=17, 0", "44÷8=5, 4"),
    @("56÷8=7, 0", "34÷9=3, 7"),
    @("47÷7=6, 5", "98÷8=12, 2"),
    @("92÷8=11, 4", "69÷3=23, 0"),
    @("35÷8=4, 3", "98÷8=12, 2"),
    @("37÷7=5, 2", "36÷4=9, 0"),
    @("37÷4=9, 1", "13÷4=3, 1"),
    @("28÷4=7, 0", "30÷7=4, 2"),
    @("73÷8=9, 1", "89÷9=9, 8"),
    @("44÷7=6, 2", "54÷8=6, 6"),
    @("54÷9=6, 0", "41÷4=10, 1"),
    @("39÷2=19, 1", "52÷4=13, 0"),
    @("23÷3=7, 2", "32÷6=5, 2"),
    @("92÷6=15, 2", "42÷9=4, 6"),
    @("53÷9=5, 8", "44÷6=7, 2"),
    @("56÷5=11, 1", "15÷9=1, 6"),
    @("42÷3=14, 0", "81÷7=11, 4"),
    @("10÷3=3, 1", "32÷8=4, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
